$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 149, pushing existing rows 149..215 down to 150..216
$ws.Rows.Item(149).Insert()

# Populate the newly inserted row 149 with the new data record
$ws.Cells.Item(149, 1).Value = 3
$ws.Cells.Item(149, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(149, 3).Value = "Coquimbo"
$ws.Cells.Item(149, 4).Value = 44489
$ws.Cells.Item(149, 5).Value = 5
$ws.Cells.Item(149, 6).Value = 100112043
$ws.Cells.Item(149, 7).Value = "Pepino ensalada"
$ws.Cells.Item(149, 8).Value = "Sin especificar"
$ws.Cells.Item(149, 9).Value = "Primera"
$ws.Cells.Item(149, 10).Value = 115
$ws.Cells.Item(149, 11).Value = 8000
$ws.Cells.Item(149, 12).Value = 8500
$ws.Cells.Item(149, 13).Value = 8261
$ws.Cells.Item(149, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(149, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(149, 16).Value = 118
$ws.Cells.Item(149, 17).Value = 70
$ws.Cells.Item(149, 18).Value = "Hortaliza"
